$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Scope 1 and 2 and 3 and 4 GHG emissions" field row (row 20);
# all subsequent rows shift up by one, and the last row (116) disappears.
$ws.Rows.Item(20).Delete() | Out-Null

# Re-establish the AutoFilter over the new, smaller range (the delete above
# does not automatically shrink the AutoFilter range).
$ws.AutoFilterMode = $false
$ws.Range("A1:M115").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$filterDbName = $wb.Names.Item("Framework Data Model!_FilterDatabase")
$filterDbName.RefersTo = "='Framework Data Model'!`$A`$1:`$M`$115"

# Update the visible selection to match the post-edit state.
$ws.Range("A2:A115").Select() | Out-Null
